$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts existing rows 4+ down by one,
# carrying the formatting of the row above into the new row).
$ws.Rows.Item(4).Insert()

# Row 3 ("test LED blinking with CRC32 ..." objective) gets the new
# combined style: colored font (like D1/D2) + wrap text + vertical
# center. Start from D1's format (colored font) then turn wrap on.
$ws.Range("D1").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").WrapText = $true

# New row 4: "Establish Acknowledgement based communicaton" task, same
# style as the new D3.
$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value = "Establish Acknowledgement based communicaton"

# Give row 3 and row 4 a "DONE" marker cell in column E, matching
# the style already used by E1/E2.
$ws.Range("E1").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").Value = "DONE"

$ws.Range("E1").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value = "DONE"

$excel.CutCopyMode = $false

# Update the selection to match the committed workbook state.
$ws.Range("E7").Select()
